$d = $word.ActiveDocument

# Helper: insert a block of new text right after $pos (collapsed range),
# then re-apply formatting on each individual piece (bold on/off) so the
# engine keeps the pieces as separate <w:r> runs instead of merging them
# back into one run.
function Insert-SplitRuns($pos, [string[]]$pieces) {
    $joined = [string]::Join('', $pieces)
    $r = $d.Range($pos, $pos)
    $r.InsertAfter($joined)

    $cursor = $pos
    foreach ($piece in $pieces) {
        $pStart = $cursor
        $pEnd = $cursor + $piece.Length
        $rb = $d.Range($pStart, $pEnd)
        $rb.Bold = 1
        $rb.Bold = 0
        $cursor = $pEnd
    }
    return $cursor
}

# --- Step 1: "One instance...project state " -> "...project state based on " ---
$text = $d.Content.Text
$oldRun1 = "One instance would be to add a pivot table of project state "
$idx1 = $text.IndexOf($oldRun1)
$rng1 = $d.Range($idx1, $idx1 + $oldRun1.Length)
$rng1.Text = "One instance would be to add a pivot table of project state based on "

# --- Steps 2-4: insert 3 new runs after run 1 ---
$text = $d.Content.Text
$anchor1 = "One instance would be to add a pivot table of project state based on "
$pos = $text.IndexOf($anchor1) + $anchor1.Length
$pos = Insert-SplitRuns $pos @(
    "the number of backers for the project",
    ", with filters for category. ",
    "However, before doing "
)

# --- Step 5: old run 2 (keeps the lastRenderedPageBreak) text swap ---
$text = $d.Content.Text
$oldRun2 = "based on average backer donation, with filters for year and category. "
$idx2 = $text.IndexOf($oldRun2)
$rng2 = $d.Range($idx2, $idx2 + $oldRun2.Length)
$rng2.Text = "this, an additional column would need to be added to the initial data tab to group the backer data into specified ranges. With that completed, a bar chart could be created to visualize the data to find "

# --- Steps 6-8: insert 3 new runs after (new) run 5 ---
$text = $d.Content.Text
$anchor2 = "With that completed, a bar chart could be created to visualize the data to find "
$pos2 = $text.IndexOf($anchor2) + $anchor2.Length
$pos2 = Insert-SplitRuns $pos2 @(
    "whether a project with more backers was more successful than one with fewer backers",
    ". Another possible option would be to include a line chart that compares the monetary ranges of the",
    " average donation to the percentage of project state. This would be similar to the already performed goal range analyses but would help identify a potential meaningful relationship. Finally, for each of the individual project states, a pie chart could "
)

# --- Step 9: old run 3 "This could be visualized through the use of a " -> new run 9 ---
$text = $d.Content.Text
$oldRun3 = "This could be visualized through the use of a "
$idx3 = $text.IndexOf($oldRun3)
$rng3 = $d.Range($idx3, $idx3 + $oldRun3.Length)
$rng3.Text = "be used to show the distribution of categories and sub-categories included in two of the already present pivot tables. These pie charts would provide a different perspective on the data compared to the stacked bar charts currently being used, as the latter can be a little difficult to view at times."
